$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51, shifting existing rows 51-189 down to 52-190.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly data point.
$ws.Range("A51").Value = 9
$ws.Range("B51").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C51").Value = "Metropolitana"
$ws.Range("D51").Value = 44525
$ws.Range("E51").Value = 13
$ws.Range("F51").Value = 300000001
$ws.Range("G51").Value = "Rabanito"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 7900
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = 2747
$ws.Range("N51").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O51").Value = "Provincia de Chacabuco"
$ws.Range("P51").Value = 27
$ws.Range("Q51").Value = 100
$ws.Range("R51").Value = "Hortaliza"
